# Update the "想去人数" (want-to-go count) column F values across the
# workbook's sheets, matching the public-repo commit that regenerated the
# data for gh-pages.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 602
$ws.Range("F3").Value = 273
$ws.Range("F8").Value = 175
$ws.Range("F10").Value = 243
$ws.Range("F11").Value = 6894
$ws.Range("F12").Value = 66
$ws.Range("F13").Value = 57
$ws.Range("F14").Value = 529
$ws.Range("F16").Value = 554
$ws.Range("F17").Value = 374
$ws.Range("F21").Value = 724
$ws.Range("F22").Value = 183
$ws.Range("F27").Value = 1895
$ws.Range("F28").Value = 529

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 278

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 290

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 290
$ws.Range("F3").Value = 602
$ws.Range("F4").Value = 273
$ws.Range("F10").Value = 175
$ws.Range("F12").Value = 243
$ws.Range("F13").Value = 6894
$ws.Range("F14").Value = 66
$ws.Range("F15").Value = 57
$ws.Range("F17").Value = 529
$ws.Range("F19").Value = 554
$ws.Range("F20").Value = 374
$ws.Range("F25").Value = 278
$ws.Range("F28").Value = 724
$ws.Range("F32").Value = 183
$ws.Range("F37").Value = 1895
$ws.Range("F38").Value = 529
